$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update Version
$ws.Range("B3").Value = "0.1.7"

# Update Status
$ws.Range("B6").Value = "draft"

# Update Date
$ws.Range("B8").Value = "2024-11-22T12:33:30-06:00"

# Update first Contact row value (was "No display for ContactDetail")
$ws.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"

# Second Contact row (row 11) gets a new value
$ws.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

# Insert a new row after row 11 for Jurisdiction (pushes old row 12.. down by one),
# copying row 12's original formatting (border/alignment/wrap) into the new row
$ws.Rows.Item(12).Copy()
$ws.Rows.Item(12).Insert()
$ws.Range("A12").Value = "Jurisdiction"
$ws.Range("B12").Value = ""
